$wb = $excel.ActiveWorkbook

# --- Content change: axis "geo" renamed to "country" ---
# The shared string is used by both the "pop" sheet (A1) and the
# "__axes__" sheet (A1); update both occurrences.
$wsPop = $wb.Worksheets.Item("pop")
$wsAxes = $wb.Worksheets.Item("__axes__")
$wsGroups = $wb.Worksheets.Item("__groups__")

$wsPop.Range("A1").Value = "country"
$wsAxes.Range("A1").Value = "country"

# --- View state: clear the stale selection on the "pop" sheet ---
$wsPop.Range("A1").Select()

# --- View state: make "__groups__" the active / selected sheet ---
$wsGroups.Activate()
